$wb = $excel.ActiveWorkbook

# --- TDD sheet: update marks for Quiz (row 8) and Assinment (row 9) ---
$tdd = $wb.Worksheets.Item("TDD")
[void]$tdd.Activate()

$tdd.Range("D8").Value = 20
$tdd.Range("G8").Value = 8
$tdd.Range("H8").Value = 9

$tdd.Range("D9").Value = 20
$tdd.Range("G9").Value = 10
$tdd.Range("H9").Value = 10

# --- SQE sheet: update marks for Project (row 6) and Assinment (row 8) ---
$sqe = $wb.Worksheets.Item("SQE")
[void]$sqe.Activate()

$sqe.Range("D6").Value = 20
$sqe.Range("G6").Value = 16

$sqe.Range("D8").Value = 60
$sqe.Range("H8").Value = 18
$sqe.Range("I8").Value = 20

[void]$sqe.Range("E13").Select()

# --- HCI sheet: selection moves, no longer the active tab ---
$hci = $wb.Worksheets.Item("HCI")
[void]$hci.Activate()
[void]$hci.Range("C10").Select()

# --- Finally, TDD becomes the active sheet/tab with its own new selection ---
[void]$tdd.Activate()
[void]$tdd.Range("E8").Select()
